$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "63.776.75"
$ws.Range("E2").Value = "  +2.95%  "

$ws.Range("D3").Value = "2.545.02"
$ws.Range("E3").Value = "  +5.66%  "

$ws.Range("E4").Value = "  -0.07%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "574.16"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +2.60%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "147.91"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +6.92%  "

$ws.Range("E7").Value = "  -0.01%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.589"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.67%  "

$ws.Range("D9").Value = "2.544.53"
$ws.Range("E9").Value = "  +5.73%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.107"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +2.51%  "

$ws.Range("E11").Value = "  +0.83%  "

$ws.Range("E12").Value = "  +1.95%  "

$ws.Range("E13").Value = "  +2.85%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "28.16"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +9.62%  "

$ws.Range("D15").Value = "2.999.03"
$ws.Range("E15").Value = "  +5.79%  "

$ws.Range("D16").Value = "63.571.91"
$ws.Range("E16").Value = "  +2.69%  "

$ws.Range("E17").Value = "  +3.13%  "

$ws.Range("D18").Value = "2.543.25"
$ws.Range("E18").Value = "  +5.09%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "11.56"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +4.71%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "341.51"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.65%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "4.36"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +3.32%  "

$ws.Range("E22").Value = "  +0.58%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "1.00"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.03%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "66.21"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.88%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.170"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.93%  "

$ws.Range("E26").Value = "  +3.61%  "

$ws.Range("E27").Value = "  +0.05%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "8.33"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.05%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.43"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +3.46%  "

$ws.Range("D30").Value = "0.0₃0835"
$ws.Range("E30").Value = "  +7.09%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.93"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +9.05%  "

$ws.Range("E32").Value = "  +3.45%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "176.79"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +3.48%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.60"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +14.10%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "419.46"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +10.61%  "

$ws.Range("E36").Value = "  +2.46%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "19.10"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +3.12%  "

$ws.Range("E38").Value = "  -2.17%  "

$ws.Range("E39").Value = "  +0.01%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.76"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +5.55%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.999"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.05%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "40.49"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +3.46%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "152.78"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +5.59%  "

$ws.Range("E44").Value = "  +3.91%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "20.94"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +1.54%  "

$ws.Range("E46").Value = "  +4.23%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0532"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.75%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0967"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.78%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "18.81"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +4.99%  "

$ws.Range("E50").Value = "  +5.40%  "

$ws.Range("D51").Value = "0.0₆0235"
$ws.Range("E51").Value = "  +8.53%  "
